$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (price + 1h volume change) per the commit diff.
# Force text format on these cells so numeric-looking strings (e.g. "65.608.32",
# "29.00", "0.0000198") are preserved exactly instead of being parsed as numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.608.32"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.54%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.666.25"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.74%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.79"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.43"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.88%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.615"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +4.64%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.67%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.31%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.81%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "29.00"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000198"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -4.00%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.142.34"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.91%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.449.54"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.48%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.669.03"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.27%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.65"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.30%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.73%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.39%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "350.94"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.07"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.44%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.73"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.52%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +3.28%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -3.17%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -4.26%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.24%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.13"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.89%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "529.91"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.32%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.53%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.60%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.27%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.58%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.59"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.26%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "156.20"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.48%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.93"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -3.00%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "162.25"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.55%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.13%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.29%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0609"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.41%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.57"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -4.10%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.33%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.82%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₆0253"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +6.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0986"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.58%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.89"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -3.53%  "
